$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) First three rows: 100 -> 0M, 0 -> 0M, 11 -> 0M
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# 2) Insert 10 new single-value rows right after the (former) third row,
#    i.e. immediately before what is currently row 4.
$values = @("12", "0.00003", "0.00010", "0.00006", "0.00000", "0.00010", "0.00010", "0.00010", "0.00051", "100.0")
$idx = 4
foreach ($v in $values) {
    $ref = $t.Rows.Item($idx)
    $newRow = $t.Rows.Add($ref)
    $newRow.Cells.Item(1).Range.Text = $v
    $idx = $idx + 1
}

# 3) Collapse the three tab-separated summary rows near the end down to a
#    single value each (these are now rows 44, 45 and 46).
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "11"
